# cycle 2 round 3
# - update cached date field text (slide master + all 3 layouts)
# - merge the split "NOC " / "Rating" runs into a single "NOC Rating" run
#   on slides 1, 4, 5, 7
# - add a new "Pain Level" textbox on slide 1
# - add new "Anxiety Level" textboxes on slides 4 and 7

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text: 9/14/2012 -> 11/8/2012
#    (slide master "Date Placeholder" + the 3 slide layouts' own copy)
# ---------------------------------------------------------------------
$master = $p.SlideMaster

foreach ($sh in $master.Shapes) {
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "11/8/2012"
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    foreach ($sh in $layout.Shapes) {
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "11/8/2012"
        }
    }
}

# ---------------------------------------------------------------------
# 2. Merge "NOC " + "Rating" runs into a single "NOC Rating" run.
#    Setting identical text is a no-op in this host, so nudge through a
#    throwaway value first to force the runs to actually collapse.
# ---------------------------------------------------------------------
function Set-MergedText($shape, $text) {
    $shape.TextFrame.TextRange.Text = "~tmp~"
    $shape.TextFrame.TextRange.Text = $text
}

$nocShapeIndex = @{1 = 10; 4 = 5; 5 = 3; 7 = 4}
foreach ($slideIdx in $nocShapeIndex.Keys) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item($nocShapeIndex[$slideIdx])
    Set-MergedText $shape "NOC Rating"
}

# ---------------------------------------------------------------------
# 3. New textboxes
# ---------------------------------------------------------------------
function Add-LabelTextbox($slide, $name, $text) {
    $left = 3138174 / 12700.0
    $top = 2244546 / 12700.0
    $width = 3150870 / 12700.0
    $height = 246221 / 12700.0

    $tb = $slide.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $tb.Name = $name

    $tb.Fill.Visible = $false

    $tf = $tb.TextFrame
    $tf.WordWrap = $true
    $tf.AutoSize = 1

    $tr = $tf.TextRange
    $tr.Text = $text
    $tr.ParagraphFormat.Alignment = 2

    $font = $tr.Font
    $font.Name = "Tahoma"
    $font.NameFarEast = "Tahoma"
    $font.NameComplexScript = "Tahoma"
    $font.Size = 10
    $font.Bold = $true

    return $tb
}

# Slide 1: "Pain Level" textbox appended after the existing shapes
$slide1 = $p.Slides.Item(1)
Add-LabelTextbox $slide1 "TextBox 21" "Pain Level" | Out-Null

# Slide 4: "Anxiety Level" textbox appended after the existing shapes
$slide4 = $p.Slides.Item(4)
Add-LabelTextbox $slide4 "TextBox 10" "Anxiety Level" | Out-Null

# Slide 7: "Anxiety Level" textbox appended after the existing shapes
$slide7 = $p.Slides.Item(7)
Add-LabelTextbox $slide7 "TextBox 12" "Anxiety Level" | Out-Null
